$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 67.2
$ws.Range("I5").Value = 67.2
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 67.2
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 47.8
$ws.Range("N5").Value = ""
$ws.Range("H11").Value = 84.75
$ws.Range("I11").Value = 84.75
$ws.Range("K11").Value = 84.75
$ws.Range("M11").Value = 55.25
$ws.Range("H18").Value = 1994.5483
$ws.Range("I18").Value = 1920.5238
$ws.Range("K18").Value = 1920.5238
$ws.Range("M18").Value = -1636.5238
$ws.Range("H40").Value = 2333.2222
$ws.Range("I40").Value = 1999.75
$ws.Range("J40").Value = 2600
$ws.Range("K40").Value = 1999.75
$ws.Range("L40").Value = 2600
$ws.Range("M40").Value = -1824.75
$ws.Range("N40").Value = -2950
$ws.Range("H69").Value = 4815
$ws.Range("J69").Value = 4815
$ws.Range("L69").Value = 14445
$ws.Range("N69").Value = -16193
$ws.Range("H72").Value = 4815
$ws.Range("J72").Value = 4815
$ws.Range("L72").Value = 43335
$ws.Range("N72").Value = -52071
$ws.Range("H116").Value = 4749.1665
$ws.Range("I116").Value = 4749.1665
$ws.Range("K116").Value = 4749.1665
$ws.Range("M116").Value = -1307.1665
$ws.Range("H138").Value = 12074.154
$ws.Range("I138").Value = 8794.200000000001
$ws.Range("J138").Value = 12855.096
$ws.Range("K138").Value = 26382.6
$ws.Range("L138").Value = 38565.288
$ws.Range("M138").Value = -21242.6
$ws.Range("N138").Value = -48845.288

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 30005
$ws.Range("J10").Value = 30005
$ws.Range("L10").Value = 30005
$ws.Range("N10").Value = -30345
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").Value = ""
$ws.Range("H12").Value = 3000000
$ws.Range("I12").Value = 3000000
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 3000000
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -2999827
$ws.Range("N12").Value = ""
$ws.Range("H16").Value = 5000
$ws.Range("I16").Value = 5000
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 5000
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -4713
$ws.Range("N16").Value = ""
$ws.Range("H61").Value = 2650.875
$ws.Range("I61").Value = 2751
$ws.Range("K61").Value = 2751
$ws.Range("M61").Value = -2539
$ws.Range("H97").Value = 2548
$ws.Range("I97").Value = 2378.8333
$ws.Range("J97").Value = 3055.5
$ws.Range("K97").Value = 2378.8333
$ws.Range("L97").Value = 3055.5
$ws.Range("M97").Value = -1882.8333
$ws.Range("N97").Value = -4047.5
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").Value = ""
$ws.Range("H136").Value = 2650.875
$ws.Range("I136").Value = 2751
$ws.Range("K136").Value = 8253
$ws.Range("M136").Value = -5703

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 60000
$ws.Range("J61").Value = 60000
$ws.Range("L61").Value = 60000
$ws.Range("N61").Value = -60626
$ws.Range("H81").Value = 40000
$ws.Range("J81").Value = 40000
$ws.Range("L81").Value = 40000
$ws.Range("N81").Value = -42122
$ws.Range("H84").Value = 40000
$ws.Range("J84").Value = 40000
$ws.Range("L84").Value = 120000
$ws.Range("N84").Value = -130608

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 5067
$ws.Range("I12").Value = 3451.6667
$ws.Range("J12").Value = 7490
$ws.Range("K12").Value = 3451.6667
$ws.Range("L12").Value = 7490
$ws.Range("M12").Value = -3281.6667
$ws.Range("N12").Value = -7830
$ws.Range("H31").Value = 2734.9412
$ws.Range("I31").Value = 3600.9092
$ws.Range("K31").Value = 3600.9092
$ws.Range("M31").Value = -3305.9092
$ws.Range("H34").Value = 2734.9412
$ws.Range("I34").Value = 3600.9092
$ws.Range("K34").Value = 3600.9092
$ws.Range("M34").Value = -3398.9092
$ws.Range("H99").Value = 11156.689
$ws.Range("I99").Value = 6930.2
$ws.Range("K99").Value = 6930.2
$ws.Range("M99").Value = -5432.2
$ws.Range("H126").Value = 11156.689
$ws.Range("I126").Value = 6930.2
$ws.Range("K126").Value = 20790.6
$ws.Range("M126").Value = -18320.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 62530.75
$ws.Range("I2").Value = 111143.22
$ws.Range("J2").Value = 29
$ws.Range("K2").Value = 666859.3200000001
$ws.Range("L2").Value = 174
$ws.Range("M2").Value = -666746.3200000001
$ws.Range("N2").Value = -400
$ws.Range("H23").Value = 189.4
$ws.Range("I23").Value = 182.33333
$ws.Range("K23").Value = 546.99999
$ws.Range("M23").Value = -311.99999
$ws.Range("H39").Value = 3249.75
$ws.Range("I39").Value = 2000
$ws.Range("J39").Value = 3666.3333
$ws.Range("K39").Value = 6000
$ws.Range("L39").Value = 10998.9999
$ws.Range("M39").Value = -5706
$ws.Range("N39").Value = -11586.9999
$ws.Range("H46").Value = 3333916.8
$ws.Range("J46").Value = 5000500
$ws.Range("L46").Value = 15001500
$ws.Range("N46").Value = -15001682
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").Value = ""
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").Value = ""
$ws.Range("H113").Value = 2123.5715
$ws.Range("J113").Value = 2533.3333
$ws.Range("L113").Value = 7599.999899999999
$ws.Range("N113").Value = -11939.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 3076.25
$ws.Range("I14").Value = 3833.3333
$ws.Range("J14").Value = 805
$ws.Range("K14").Value = 3833.3333
$ws.Range("L14").Value = 805
$ws.Range("M14").Value = -3665.3333
$ws.Range("N14").Value = -1141
$ws.Range("H122").Value = 33300.094
$ws.Range("I122").Value = 1978.6818
$ws.Range("J122").Value = 102207.2
$ws.Range("K122").Value = 5936.0454
$ws.Range("L122").Value = 306621.6
$ws.Range("M122").Value = -3486.0454
$ws.Range("N122").Value = -311521.6
$ws.Range("H132").Value = 2768.425
$ws.Range("I132").Value = 2570
$ws.Range("J132").Value = 2987.7368
$ws.Range("K132").Value = 7710
$ws.Range("L132").Value = 8963.2104
$ws.Range("M132").Value = -5180
$ws.Range("N132").Value = -14023.2104

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 10000
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 10000
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 10000
$ws.Range("M10").Value = ""
$ws.Range("N10").Value = -10280
$ws.Range("H18").Value = 37502
$ws.Range("I18").Value = 37502
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 37502
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -37330
$ws.Range("N18").Value = ""
$ws.Range("H22").Value = 6581.1934
$ws.Range("I22").Value = 5483.7856
$ws.Range("J22").Value = 7484.9414
$ws.Range("K22").Value = 5483.7856
$ws.Range("L22").Value = 7484.9414
$ws.Range("M22").Value = -5188.7856
$ws.Range("N22").Value = -8074.9414
$ws.Range("H24").Value = 2748.5
$ws.Range("J24").Value = 2748.5
$ws.Range("L24").Value = 2748.5
$ws.Range("N24").Value = -3434.5
$ws.Range("H27").Value = 6581.1934
$ws.Range("I27").Value = 5483.7856
$ws.Range("J27").Value = 7484.9414
$ws.Range("K27").Value = 5483.7856
$ws.Range("L27").Value = 7484.9414
$ws.Range("M27").Value = -5376.7856
$ws.Range("N27").Value = -7698.9414
$ws.Range("H40").Value = 3721.75
$ws.Range("I40").Value = 2962.3333
$ws.Range("K40").Value = 2962.3333
$ws.Range("M40").Value = -2826.3333
$ws.Range("H122").Value = 4641.8
$ws.Range("I122").Value = 4552.25
$ws.Range("K122").Value = 13656.75
$ws.Range("M122").Value = -11206.75
$ws.Range("H132").Value = 5668.5
$ws.Range("I132").Value = 5075.1665
$ws.Range("J132").Value = 6024.5
$ws.Range("K132").Value = 15225.4995
$ws.Range("L132").Value = 18073.5
$ws.Range("M132").Value = -12695.4995
$ws.Range("N132").Value = -23133.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 999999
$ws.Range("I10").Value = 999999
$ws.Range("K10").Value = 999999
$ws.Range("M10").Value = -999830
$ws.Range("H14").Value = 500450
$ws.Range("I14").Value = 500450
$ws.Range("K14").Value = 500450
$ws.Range("M14").Value = -500282
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").Value = ""
$ws.Range("H132").Value = 119453.625
$ws.Range("I132").Value = 237008.5
$ws.Range("K132").Value = 711025.5
$ws.Range("M132").Value = -708495.5
